$d = $word.ActiveDocument

$replacements = @(
    @("67÷6=", "20÷9="),
    @("93÷7=", "85÷8="),
    @("48÷8=", "71÷9="),
    @("14÷3=", "36÷5="),
    @("10÷8=", "60÷8="),
    @("98÷5=", "96÷5="),
    @("96÷3=", "41÷5="),
    @("37÷2=", "67÷4="),
    @("25÷8=", "52÷2="),
    @("53÷7=", "44÷7="),
    @("63÷9=", "93÷5="),
    @("19÷9=", "14÷6="),
    @("10÷3=", "49÷3="),
    @("21÷3=", "77÷8="),
    @("81÷7=", "46÷4="),
    @("46÷8=", "25÷8="),
    @("94÷9=", "19÷6="),
    @("56÷9=", "75÷8="),
    @("94÷4=", "61÷8="),
    @("99÷3=", "18÷2="),
    @("23÷5=", "89÷4="),
    @("58÷5=", "85÷3="),
    @("93÷9=", "94÷2="),
    @("82÷5=", "80÷5="),
    @("47÷5=", "62÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
